# "added family promise census info"
# Refresh the exit-destination percentages with the latest Family Promise
# census numbers. This both updates existing rows' labels/values and
# inserts several new exit-destination categories, growing the table
# from A1:B8 to A1:B14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "Client refused",
    "Emergency shelter, including hotel or motel paid for with emergency shelter voucher, or RHY-funded Host Home shelter",
    "No exit interview completed",
    "Rental by client, no ongoing housing subsidy",
    "Rental by client with RRH or equivalent subsidy",
    "Staying or living with friends, temporary tenure (e.g., room, apartment or house)",
    "Staying or living with family, permanent tenure",
    "Transitional Housing for homeless persons (including homeless youth)",
    "Hotel or Motel paid for without Emergency Shelter Voucher",
    "Staying or living with family, temporary tenure (e.g., room, apartment or house)",
    "Other",
    "Substance Abuse Treatment or Detox Center",
    "Rental by client, other ongoing housing subsidy"
)

$values = @(
    0.3981900452488688,
    0.16289592760181,
    0.09954751131221719,
    0.08144796380090498,
    0.07692307692307693,
    0.04977375565610859,
    0.04072398190045249,
    0.03619909502262444,
    0.01809954751131222,
    0.01357466063348416,
    0.009049773755656109,
    0.009049773755656109,
    0.004524886877828055
)

# Pre-existing rows already carry the right formatting (style "s=1" on
# column A) - only the new rows 9-14 need that style cloned. Use the
# already-formatted A2 cell as the format source, copy/paste-special so
# the destination reuses the same style entry instead of growing a
# near-duplicate one.
$ws.Range("A2").Copy() | Out-Null
for ($row = 9; $row -le 14; $row++) {
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
